$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$vb = [char]11

$t.Cell(1, 1).Range.Text = "91 x 19" + $vb + "  1    9" + $vb + "  ----" + $vb + "9|    |" + $vb + "1|    |"
$t.Cell(1, 2).Range.Text = "59 x 59" + $vb + "  5    9" + $vb + "  ----" + $vb + "5|    |" + $vb + "9|    |"
$t.Cell(1, 3).Range.Text = "39 x 35" + $vb + "  3    5" + $vb + "  ----" + $vb + "3|    |" + $vb + "9|    |"
$t.Cell(2, 1).Range.Text = "89 x 15" + $vb + "  1    5" + $vb + "  ----" + $vb + "8|    |" + $vb + "9|    |"
$t.Cell(2, 2).Range.Text = "78 x 85" + $vb + "  8    5" + $vb + "  ----" + $vb + "7|    |" + $vb + "8|    |"
$t.Cell(2, 3).Range.Text = "83 x 21" + $vb + "  2    1" + $vb + "  ----" + $vb + "8|    |" + $vb + "3|    |"
$t.Cell(3, 1).Range.Text = "52 x 93" + $vb + "  9    3" + $vb + "  ----" + $vb + "5|    |" + $vb + "2|    |"
$t.Cell(3, 2).Range.Text = "76 x 11" + $vb + "  1    1" + $vb + "  ----" + $vb + "7|    |" + $vb + "6|    |"
$t.Cell(3, 3).Range.Text = "71 x 31" + $vb + "  3    1" + $vb + "  ----" + $vb + "7|    |" + $vb + "1|    |"
$t.Cell(4, 1).Range.Text = "34 x 34" + $vb + "  3    4" + $vb + "  ----" + $vb + "3|    |" + $vb + "4|    |"
$t.Cell(4, 2).Range.Text = "73 x 40" + $vb + "  4    0" + $vb + "  ----" + $vb + "7|    |" + $vb + "3|    |"
$t.Cell(4, 3).Range.Text = "45 x 33" + $vb + "  3    3" + $vb + "  ----" + $vb + "4|    |" + $vb + "5|    |"
$t.Cell(5, 1).Range.Text = "51 x 10" + $vb + "  1    0" + $vb + "  ----" + $vb + "5|    |" + $vb + "1|    |"
$t.Cell(5, 2).Range.Text = "97 x 16" + $vb + "  1    6" + $vb + "  ----" + $vb + "9|    |" + $vb + "7|    |"
$t.Cell(5, 3).Range.Text = "57 x 34" + $vb + "  3    4" + $vb + "  ----" + $vb + "5|    |" + $vb + "7|    |"
